$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the bold/border/centered style previously applied to A1, and drop
# the second data row (A2) whose shared-string content is being merged
# into A1 as reformatted JSON text.
$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()

$ws.Range("A1").Value = @"
questions = [
    {
        "title": "You are a network consultant. A small office has approached you for advice. They recently expanded their team from two people to 25 people, all within a single office. They need an improved networking system to support their operations efficiently. What would you recommend for establishing a networking solution within their small office?",
        "ques_type": 2,
        "options": [
            "Local area network (LAN)",
            "Metropolitan area network (MAN)",
            "Wide area network (WAN)",
            "Personal area network (PAN)"
        ],
        "score": "Local area network (LAN)"
    },
    {
        "title": "You are a network engineer. A client has asked you to set up a network within an office with low latency. You should be able to add or remove devices from the network without modifying its architecture. The client wants you to determine the best method to achieve this.Which network should you recommend?",
        "ques_type": 2,
        "options": [
            "Bluetooth-based network",
            "Broadband-based wireless network",
            "Ethernet wired network",
            "Hotspot network"
        ],
        "score": "Ethernet wired network"
    },
    {
        "title": "You are a network engineer. Users In a Transmission Control Protocol/Internet Protocol (TCP/IP) network report connectivity issues when downloading a large file from a remote server. The users suspect that the connection is unreliable. Which statement correctly describes the problem?",
        "ques_type": 2,
        "options": [
            "A misconfiguration of Domain Name System (DNS) settings",
            "A problem with a connectionless (stateless) protocol",
            "An error in the routing tables",
            "A failure in a connection-oriented (stateful) protocol"
        ],
        "score": "A failure in a connection-oriented (stateful) protocol"
    },
    {
        "title": "You are a network consultant and have advised your client to implement cryptography within the network. Your client has asked for the reason behind this recommendation.What should be your answer to the client?",
        "ques_type": 2,
        "options": [
            "It ensures high speed and low latency.",
            "It facilitates rapid data transmission without encryption.",
            "It monitors and analyzes network traffic for vulnerabilities.",
            "It protects data confidentiality, integrity, and authenticity."
        ],
        "score": "It protects data confidentiality, integrity, and authenticity."
    }
]
"@

# Re-entering the value introduced embedded newlines, which Excel uses to
# auto-grow the row height; put the row back to its default auto height.
$ws.Rows(1).AutoFit()
